$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.854886412620544
$ws.Range("B1").Value = 3.140937089920044
$ws.Range("C1").Value = 2.705324172973633
$ws.Range("D1").Value = 2.202339172363281
$ws.Range("E1").Value = 1.45538318157196
